# ----------------------------------------------------------------------
# Applies the "Added some client side display for upvoting" edit:
#   1. Splits "Upvoting and downvoting of questions and answers" into
#      separate runs with spell-check proofErr markers around the two
#      non-dictionary words "Upvoting" / "downvoting".
#   2. Expands "Push notification system " into
#      "Push notification system / Fix button things in notification windows".
#   3. Inserts a new "Break up index.js page" bullet and a new
#      "Add commenting" bullet after "Clean up code", moving the
#      "_GoBack" bookmark from the old "Add commenting" paragraph onto
#      the "Nice to have" paragraph that follows it.
# ----------------------------------------------------------------------

$d = $word.ActiveDocument

$rPr = '<w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="373E4D"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:color="auto" w:fill="F6F7F8"/></w:rPr>'

function New-PackageXml([string]$body) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ------------------------------------------------------------------
# 1. "Upvoting and downvoting of questions and answers" -> 4 runs
#    with proofErr spell-check markers around "Upvoting" / "downvoting"
# ------------------------------------------------------------------
$found = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Upvoting and downvoting of questions and answers") {
        $found = $p
        break
    }
}

$full = $found.Range
$textRange = $d.Range($full.Start, $full.End - 1)

$upvoteBody = ''
$upvoteBody += '<w:proofErr w:type="spellStart"/>'
$upvoteBody += '<w:r>' + $rPr + '<w:t>Upvoting</w:t></w:r>'
$upvoteBody += '<w:proofErr w:type="spellEnd"/>'
$upvoteBody += '<w:r>' + $rPr + '<w:t xml:space="preserve"> and </w:t></w:r>'
$upvoteBody += '<w:proofErr w:type="spellStart"/>'
$upvoteBody += '<w:r>' + $rPr + '<w:t>downvoting</w:t></w:r>'
$upvoteBody += '<w:proofErr w:type="spellEnd"/>'
$upvoteBody += '<w:r>' + $rPr + '<w:t xml:space="preserve"> of questions and answers</w:t></w:r>'

# The target range sits inside an existing paragraph (it excludes the
# paragraph mark), so the replacement fragment must still be wrapped in
# its own <w:p> for InsertXML to treat it as "this paragraph's content"
# rather than one implicit (and, here, lost) paragraph per top-level node.
$textRange.InsertXML((New-PackageXml ('<w:p>' + $upvoteBody + '</w:p>')))

# ------------------------------------------------------------------
# 2. "Push notification system " -> "Push notification system / Fix
#    button things in notification windows"
# ------------------------------------------------------------------
$pushPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Push notification system") {
        $pushPara = $p
        break
    }
}
$pushFull = $pushPara.Range
$pushTextRange = $d.Range($pushFull.Start, $pushFull.End - 1)
$pushBody = '<w:p><w:r w:rsidRPr="00DC6F71">' + $rPr + '<w:t>Push notification system / Fix button things in notification windows</w:t></w:r></w:p>'
$pushTextRange.InsertXML((New-PackageXml $pushBody))

# ------------------------------------------------------------------
# 3. Insert "Break up index.js page" and "Add commenting" bullets
#    after "Clean up code"; move the "_GoBack" bookmark from the old
#    "Add commenting" paragraph to the "Nice to have" paragraph.
# ------------------------------------------------------------------
$addCommentingPara = $null
$niceToHavePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Add commenting") {
        $addCommentingPara = $p
    }
    if ($p.Range.Text.TrimEnd() -eq "Nice to have") {
        $niceToHavePara = $p
    }
}

$rangeStart = $addCommentingPara.Range.Start
$rangeEnd = $niceToHavePara.Range.End
$replaceRange = $d.Range($rangeStart, $rangeEnd)

$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$indPPr = '<w:pPr><w:ind w:left="360"/></w:pPr>'

$breakUpPara = '<w:p>' + $listPPr + '<w:r>' + $rPr + '<w:t>Break up index.js page</w:t></w:r></w:p>'
$addCommentingNewPara = '<w:p>' + $listPPr + '<w:r>' + $rPr + '<w:t>Add commenting</w:t></w:r></w:p>'
$niceToHaveNewPara = '<w:p>' + $indPPr + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="00DC6F71">' + $rPr + '<w:t>Nice to have</w:t></w:r></w:p>'

$replaceRange.InsertXML((New-PackageXml ($breakUpPara + $addCommentingNewPara + $niceToHaveNewPara)))
